$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aquisicoes")
$ws.Name = "Folha1"
